$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3703.0527
$ws.Range("I62").Value = 2657.6667
$ws.Range("J62").Value = 5495.143
$ws.Range("K62").Value = 2657.6667
$ws.Range("L62").Value = 5495.143
$ws.Range("M62").Value = -2033.6667
$ws.Range("N62").Value = -6743.143
$ws.Range("H65").Value = 3703.0527
$ws.Range("I65").Value = 2657.6667
$ws.Range("J65").Value = 5495.143
$ws.Range("K65").Value = 13288.3335
$ws.Range("L65").Value = 27475.715
$ws.Range("M65").Value = -10168.3335
$ws.Range("N65").Value = -33715.715
$ws.Range("H92").Value = 2100.5715
$ws.Range("I92").Value = 1940.8
$ws.Range("J92").Value = 2500
$ws.Range("K92").Value = 1940.8
$ws.Range("L92").Value = 2500
$ws.Range("M92").Value = -692.8
$ws.Range("N92").Value = -4996
$ws.Range("H113").Value = 3764.0908
$ws.Range("I113").Value = 2882
$ws.Range("J113").Value = 4499.1665
$ws.Range("K113").Value = 2882
$ws.Range("L113").Value = 4499.1665
$ws.Range("M113").Value = 372
$ws.Range("N113").Value = -11007.1665
$ws.Range("H129").Value = 889.4167
$ws.Range("I129").Value = 261.75
$ws.Range("J129").Value = 1014.95
$ws.Range("K129").Value = 785.25
$ws.Range("L129").Value = 3044.85
$ws.Range("M129").Value = 4214.75
$ws.Range("N129").Value = -13044.85
$ws.Range("H137").Value = 1982.591
$ws.Range("I137").Value = 1713.3334
$ws.Range("J137").Value = 3194.25
$ws.Range("K137").Value = 5140.0002
$ws.Range("L137").Value = 9582.75
$ws.Range("M137").Value = -2590.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 341.66666
$ws.Range("I5").Value = 330
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 330
$ws.Range("L5").Value = 400
$ws.Range("M5").Value = -218
$ws.Range("N5").Value = -624
$ws.Range("H23").Value = 18572.143
$ws.Range("I23").Value = 5005
$ws.Range("J23").Value = 20833.334
$ws.Range("K23").Value = 5005
$ws.Range("L23").Value = 20833.334
$ws.Range("M23").Value = -4746
$ws.Range("N23").Value = -21351.334
$ws.Range("H61").Value = 1679.6875
$ws.Range("I61").Value = 1484.8
$ws.Range("J61").Value = 2004.5
$ws.Range("K61").Value = 1484.8
$ws.Range("L61").Value = 2004.5
$ws.Range("M61").Value = -1272.8
$ws.Range("N61").Value = -2428.5
$ws.Range("H107").Value = 29283
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 29283
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 29283
$ws.Range("N107").Value = -36963
$ws.Range("H122").Value = 1588.1052
$ws.Range("I122").Value = 1119.7273
$ws.Range("J122").Value = 2232.125
$ws.Range("K122").Value = 3359.1819
$ws.Range("L122").Value = 6696.375
$ws.Range("M122").Value = -909.1819
$ws.Range("N122").Value = -11596.375
$ws.Range("H136").Value = 1679.6875
$ws.Range("I136").Value = 1484.8
$ws.Range("J136").Value = 2004.5
$ws.Range("K136").Value = 4454.4
$ws.Range("L136").Value = 6013.5
$ws.Range("M136").Value = -1904.4
$ws.Range("N136").Value = -11113.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 341.66666
$ws.Range("I4").Value = 330
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 330
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -215
$ws.Range("N4").Value = -630
$ws.Range("H36").Value = 6767.143
$ws.Range("I36").Value = 1474
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 1474
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -940
$ws.Range("N36").Value = -21068
$ws.Range("H99").Value = 2081.8333
$ws.Range("I99").Value = 1690
$ws.Range("J99").Value = 2160.2
$ws.Range("K99").Value = 1690
$ws.Range("L99").Value = 2160.2
$ws.Range("M99").Value = -192
$ws.Range("N99").Value = -5156.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 386.6
$ws.Range("I22").Value = 156.66667
$ws.Range("J22").Value = 731.5
$ws.Range("K22").Value = 156.66667
$ws.Range("L22").Value = 731.5
$ws.Range("M22").Value = 193.33333
$ws.Range("N22").Value = -1431.5
$ws.Range("H132").Value = 3828.476
$ws.Range("I132").Value = 3268
$ws.Range("J132").Value = 4248.8335
$ws.Range("K132").Value = 9804
$ws.Range("L132").Value = 12746.5005
$ws.Range("M132").Value = -7274

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 537.375
$ws.Range("I34").Value = 375.75
$ws.Range("J34").Value = 699
$ws.Range("K34").Value = 1127.25
$ws.Range("L34").Value = 2097
$ws.Range("M34").Value = -1043.25
$ws.Range("N34").Value = -2265
$ws.Range("H124").Value = 226
$ws.Range("I124").Value = 226
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 678
$ws.Range("L124").Value = 0
$ws.Range("M124").Value = 4232
$ws.Range("H131").Value = 886.75
$ws.Range("I131").Value = 564.1429000000001
$ws.Range("J131").Value = 911.0323
$ws.Range("K131").Value = 1692.4287
$ws.Range("L131").Value = 2733.0969
$ws.Range("M131").Value = 3347.5713
$ws.Range("N131").Value = -12813.0969

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 37254.5
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 37254.5
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 37254.5
$ws.Range("N25").Value = -38312.5
$ws.Range("H51").Value = 29333.334
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 29333.334
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 29333.334
$ws.Range("N51").Value = -30351.334
$ws.Range("H97").Value = 1788.8889
$ws.Range("I97").Value = 1620
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1620
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -1124
$ws.Range("I102").Value = 1814.3
$ws.Range("J102").Value = 1964.125
$ws.Range("K102").Value = 1814.3
$ws.Range("L102").Value = 1964.125
$ws.Range("M102").Value = -192.3
$ws.Range("N102").Value = -5208.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3214.3
$ws.Range("I7").Value = 3237.5557
$ws.Range("J7").Value = 3005
$ws.Range("K7").Value = 3237.5557
$ws.Range("L7").Value = 3005
$ws.Range("M7").Value = -3125.5557
$ws.Range("H22").Value = 724.4783
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 734.6818
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 734.6818
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1324.6818
$ws.Range("H27").Value = 724.4783
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 734.6818
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 734.6818
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -948.6818
$ws.Range("H40").Value = 2050
$ws.Range("I40").Value = 1766.6666
$ws.Range("J40").Value = 2900
$ws.Range("K40").Value = 1766.6666
$ws.Range("L40").Value = 2900
$ws.Range("M40").Value = -1630.6666
$ws.Range("N40").Value = -3172
$ws.Range("H93").Value = 88050.875
$ws.Range("I93").Value = 629
$ws.Range("J93").Value = 700004
$ws.Range("K93").Value = 629
$ws.Range("L93").Value = 700004
$ws.Range("M93").Value = 619
$ws.Range("H122").Value = 4198.926
$ws.Range("I122").Value = 4724
$ws.Range("J122").Value = 2951.875
$ws.Range("K122").Value = 14172
$ws.Range("L122").Value = 8855.625
$ws.Range("M122").Value = -11722
$ws.Range("N122").Value = -13755.625
$ws.Range("H126").Value = 3214.3
$ws.Range("I126").Value = 3237.5557
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 9712.667099999999
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -7242.667099999999
$ws.Range("H136").Value = 12821894
$ws.Range("I136").Value = 15152725
$ws.Range("J136").Value = 2325
$ws.Range("K136").Value = 45458175
$ws.Range("L136").Value = 6975
$ws.Range("M136").Value = -45455625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 71429336
$ws.Range("I126").Value = 125000650
$ws.Range("J126").Value = 914.8333
$ws.Range("K126").Value = 375001950
$ws.Range("L126").Value = 2744.4999
$ws.Range("M126").Value = -374999480
$ws.Range("N126").Value = -7684.4999
$ws.Range("H136").Value = 20410114
$ws.Range("I136").Value = 32260398
$ws.Range("J136").Value = 1291.2778
$ws.Range("K136").Value = 96781194
$ws.Range("L136").Value = 3873.8334
$ws.Range("M136").Value = -96778644
$ws.Range("N136").Value = -8973.8334
